# Atualização de bases das ligas, do dia: 13-06-2024 às 19:35
#
# The source data for three pairs of matches (rows 73/74, 112/113, 159/160)
# had their Home/Away team assignment - and therefore every stat/odds column
# that goes with it - swapped. Columns A (row index), C (Div) and D (Date)
# are unaffected; only B (id) and E:AD (HomeTeam .. PL_AhUnder) need to be
# exchanged between the two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $range1 = $ws.Range("B$row1`:AD$row1")
    $range2 = $ws.Range("B$row2`:AD$row2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value = $vals2
    $range2.Value = $vals1
}

Swap-Rows 73 74
Swap-Rows 112 113
Swap-Rows 159 160
